$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows for September 2024 (mois_annee serial 45536), mirroring the
# existing July/August blocks (rows 2-6 and 7-11).
$title = "Cosy Appart - Plage à 2min - Casino à 1min"

$data = @(
    @{ Row = 12; Type = "electricite"; Charge = 45  },
    @{ Row = 13; Type = "copro";       Charge = 70  },
    @{ Row = 14; Type = "box ";        Charge = 30  },
    @{ Row = 15; Type = "credit";      Charge = 880 },
    @{ Row = 16; Type = "samantha";    Charge = 120 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.Type
    $ws.Range("B$r").Value = $item.Charge
    $ws.Range("C$r").Value = 45536
    $ws.Range("C$r").NumberFormat = "mmm-yy"
    $ws.Range("D$r").Value = $title
}

$ws.Range("I24").Select()
